$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New data columns F and G (with header strings in row 6) ---
$ws.Range("F1").Value = 0
$ws.Range("G1").Value = 0

$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

$ws.Range("F3").Value = 162
$ws.Range("G3").Value = 36

$ws.Range("F4").Value = 2077
$ws.Range("G4").Value = 677

$ws.Range("F5").Value = 39660
$ws.Range("G5").Value = 6075

$ws.Range("F6").Value = "C array"
$ws.Range("G6").Value = "C array parallel"

# --- 2. Column G width (closest value COM's pixel-snapped ColumnWidth can reach to 13.85546875) ---
$ws.Columns.Item(7).ColumnWidth = 13

# --- 3. Selection moves to G2 (mirrors the diff's <selection activeCell="G2" sqref="G2"/>) ---
$ws.Range("G2").Select()

# --- 4. Move the chart to make room for the new columns ---
$co = $ws.ChartObjects().Item(1)
$co.Left = 759.7411321973425
$co.Top = 49.12488188976378
$co.Width = 881.0625787401575
$co.Height = 453.3751181102362
